$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '26.319.42'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.690.86'
$ws.Cells.Item(3, 5).Value = '  +1.32%  '

$ws.Cells.Item(4, 5).Value = '  +0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '218.90'
$ws.Cells.Item(5, 5).Value = '  +0.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5273'
$ws.Cells.Item(6, 5).Value = '  +4.22%  '

$ws.Cells.Item(7, 5).Value = '  +0.07%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2707'
$ws.Cells.Item(8, 5).Value = '  +1.85%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06438'
$ws.Cells.Item(9, 5).Value = '  +1.66%  '

$ws.Cells.Item(10, 5).Value = '  +2.50%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07478'
$ws.Cells.Item(11, 5).Value = '  +1.44%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.718.57'
$ws.Cells.Item(12, 5).Value = '  +2.92%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.567'
$ws.Cells.Item(13, 5).Value = '  +0.68%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.5860'
$ws.Cells.Item(14, 5).Value = '  +1.20%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.000008548'
$ws.Cells.Item(15, 5).Value = '  +0.28%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.63'
$ws.Cells.Item(16, 5).Value = '  -0.26%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '26.360.36'
$ws.Cells.Item(17, 5).Value = '  +0.76%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '4.958'
$ws.Cells.Item(18, 5).Value = '  +0.51%  '

$ws.Cells.Item(19, 5).Value = '  +0.11%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.92'
$ws.Cells.Item(20, 5).Value = '  +0.85%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '190.09'
$ws.Cells.Item(21, 5).Value = '  +0.59%  '

$ws.Cells.Item(22, 5).Value = '  +0.63%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '144.52'
$ws.Cells.Item(24, 5).Value = '  +0.55%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '7.703'
$ws.Cells.Item(25, 5).Value = '  +0.27%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1235'
$ws.Cells.Item(26, 5).Value = '  +5.54%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.92'
$ws.Cells.Item(27, 5).Value = '  +1.46%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.06687'
$ws.Cells.Item(28, 5).Value = '  +15.70%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.360'
$ws.Cells.Item(29, 5).Value = '  +6.44%  '

$ws.Cells.Item(30, 5).Value = '  +0.73%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.590'
$ws.Cells.Item(31, 5).Value = '  +2.26%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.581'
$ws.Cells.Item(32, 5).Value = '  +1.52%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.673'
$ws.Cells.Item(33, 5).Value = '  +2.49%  '

$ws.Cells.Item(34, 5).Value = '  +2.16%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.6236'
$ws.Cells.Item(35, 5).Value = '  +4.40%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.394'
$ws.Cells.Item(36, 5).Value = '  +1.40%  '

$ws.Cells.Item(37, 5).Value = '  +2.26%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '6.382'
$ws.Cells.Item(38, 5).Value = '  +6.31%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.112.06'
$ws.Cells.Item(39, 5).Value = '  +3.72%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01626'
$ws.Cells.Item(40, 5).Value = '  +1.05%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8882'
$ws.Cells.Item(41, 5).Value = '  +3.18%  '

$ws.Cells.Item(42, 5).Value = '  +0.80%  '

$ws.Cells.Item(43, 5).Value = '  +1.35%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.838.65'
$ws.Cells.Item(44, 5).Value = '  +1.14%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00000000113'
$ws.Cells.Item(45, 5).Value = '  +1.99%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '57.05'
$ws.Cells.Item(46, 5).Value = '  +2.50%  '

$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.198'
$ws.Cells.Item(47, 5).Value = '  +1.36%  '

$ws.Cells.Item(48, 2).Value = 'Frax'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.010'
$ws.Cells.Item(48, 5).Value = '  +0.56%  '

$ws.Cells.Item(49, 5).Value = '  +1.73%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.4302'
$ws.Cells.Item(50, 5).Value = '  -0.01%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '6.070'
